$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the offset values for row 8 (K8, L8, M8). B8/C8 are formulas that
# depend on these and will recalculate automatically.
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 5
$ws.Range("M8").Value = 0

# Move the active selection to L9 (was D9).
$ws.Range("L9").Select()
